$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.802.87'
$ws.Range('E2').Value = '  -0.47%  '
$ws.Range('D3').Value = '1.892.03'
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7954'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -4.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.85'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9993'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3167'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  -2.48%  '
$ws.Range('E9').Value = '  -4.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07038'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +0.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08047'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7673'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +2.58%  '
$ws.Range('D13').Value = '1.882.99'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.296'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.80%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.24'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').Value = '29.800.99'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.88'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -1.25%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.930'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '243.72'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007718'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.215'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +18.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9994'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = '2.138.78'
$ws.Range('E23').Value = '  -0.68%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.9995'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1664'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +4.41%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.314'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +1.32%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '165.87'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  -1.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.68'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.055'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -1.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.394'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +1.73%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.537'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +1.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.417'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +3.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05640'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.047'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.262'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -1.22%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7395'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.002'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.645'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -2.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01906'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -0.15%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.774'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4413'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.48'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +0.82%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.806'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -2.26%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8420'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9988'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -0.08%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.50'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  +1.05%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '1.024.97'
$ws.Range('E47').Value = '  +3.46%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.870'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -1.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.914'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +1.70%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.428'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -2.30%  '
$ws.Range('D51').Value = '2.032.58'
$ws.Range('E51').Value = '  -0.89%  '
